$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at position 467 (pushes existing rows 467:557 down to 468:558)
$ws.Rows("467:467").Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A467").Value = 5
$ws.Range("B467").Value = "Macroferia Regional de Talca"
$ws.Range("C467").Value = "Maule"
$ws.Range("D467").Value = "2023-10-12"
$ws.Range("E467").Value = 7
$ws.Range("F467").Value = "Fruta"
$ws.Range("G467").Value = 100102
$ws.Range("H467").Value = "Cítricos"
$ws.Range("I467").Value = 100102004
$ws.Range("J467").Value = "Mandarina"
$ws.Range("K467").Value = "Murcott"
$ws.Range("L467").Value = "Primera"
$ws.Range("M467").Value = 250
$ws.Range("N467").Value = 8000
$ws.Range("O467").Value = 8000
$ws.Range("P467").Value = 8000
$ws.Range("Q467").Value = "$/bandeja 18 kilos"
$ws.Range("R467").Value = "Región de O'Higgins"
$ws.Range("S467").Value = 444
$ws.Range("T467").Value = 18
